$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $rng = $ws.Cells.Item($row, $col)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 2 4 '29.719.65'
Set-TextValue 2 5 '  +1.90%  '
Set-TextValue 3 4 '1.856.47'
Set-TextValue 3 5 '  +1.55%  '
Set-TextValue 4 4 '1.001'
Set-TextValue 4 5 '  +0.19%  '
Set-TextValue 5 4 '244.44'
Set-TextValue 5 5 '  +0.88%  '
Set-TextValue 6 4 '0.6392'
Set-TextValue 6 5 '  +3.20%  '
Set-TextValue 7 4 '1.001'
Set-TextValue 7 5 '  +0.10%  '
Set-TextValue 8 4 '46.76'
Set-TextValue 8 5 '  +3.00%  '
Set-TextValue 9 4 '0.3010'
Set-TextValue 9 5 '  +3.23%  '
Set-TextValue 10 4 '0.07480'
Set-TextValue 10 5 '  +1.70%  '
Set-TextValue 11 4 '24.27'
Set-TextValue 11 5 '  +5.15%  '
Set-TextValue 12 4 '0.07660'
Set-TextValue 12 5 '  -0.12%  '
Set-TextValue 13 4 '1.870.01'
Set-TextValue 13 5 '  +2.48%  '
Set-TextValue 14 4 '5.047'
Set-TextValue 14 5 '  +2.00%  '
Set-TextValue 15 4 '0.6878'
Set-TextValue 15 5 '  +3.64%  '
Set-TextValue 16 4 '83.92'
Set-TextValue 16 5 '  +2.15%  '
Set-TextValue 17 4 '0.000009524'
Set-TextValue 17 5 '  +6.92%  '
Set-TextValue 18 4 '6.056'
Set-TextValue 18 5 '  +3.83%  '
Set-TextValue 19 4 '29.757.65'
Set-TextValue 19 5 '  +2.18%  '
Set-TextValue 20 4 '2.121.17'
Set-TextValue 20 5 '  +2.71%  '
Set-TextValue 21 4 '236.55'
Set-TextValue 21 5 '  -0.76%  '
Set-TextValue 22 4 '12.63'
Set-TextValue 22 5 '  +1.51%  '
Set-TextValue 23 4 '1.001'
Set-TextValue 23 5 '  +0.13%  '
Set-TextValue 24 4 '7.423'
Set-TextValue 24 5 '  +1.07%  '
Set-TextValue 25 4 '1.003'
Set-TextValue 25 5 '  +0.20%  '
Set-TextValue 26 4 '158.21'
Set-TextValue 26 5 '  +0.20%  '
Set-TextValue 27 4 '0.1423'
Set-TextValue 27 5 '  +0.61%  '
Set-TextValue 28 4 '8.498'
Set-TextValue 28 5 '  -0.12%  '
Set-TextValue 29 5 '  +1.58%  '
Set-TextValue 30 4 '0.06145'
Set-TextValue 30 5 '  +3.86%  '
Set-TextValue 31 4 '1.492'
Set-TextValue 31 5 '  +0.39%  '
Set-TextValue 32 4 '1.266'
Set-TextValue 32 5 '  +4.76%  '
Set-TextValue 33 4 '4.151'
Set-TextValue 33 5 '  +1.75%  '
Set-TextValue 34 4 '4.094'
Set-TextValue 34 5 '  +0.63%  '
Set-TextValue 35 4 '1.882'
Set-TextValue 35 5 '  +0.95%  '
Set-TextValue 36 4 '1.169'
Set-TextValue 36 5 '  +3.00%  '
Set-TextValue 37 4 '0.7266'
Set-TextValue 37 5 '  -0.66%  '
Set-TextValue 38 5 '  +0.02%  '
Set-TextValue 39 4 '2.856'
Set-TextValue 39 5 '  +0.30%  '
Set-TextValue 40 4 '0.01781'
Set-TextValue 40 5 '  +1.86%  '
Set-TextValue 41 4 '1.208.10'
Set-TextValue 41 5 '  -0.66%  '
Set-TextValue 42 4 '0.9255'
Set-TextValue 42 5 '  +1.08%  '
Set-TextValue 43 4 '6.160'
Set-TextValue 43 5 '  -1.98%  '
Set-TextValue 44 4 '1.002'
Set-TextValue 44 5 '  +0.14%  '
Set-TextValue 45 4 '2.030.39'
Set-TextValue 45 5 '  +3.10%  '
Set-TextValue 46 4 '102.00'
Set-TextValue 46 5 '  +0.10%  '
Set-TextValue 47 4 '66.23'
Set-TextValue 47 5 '  +2.28%  '
Set-TextValue 48 4 '0.00000000123'
Set-TextValue 48 5 '  +5.30%  '
Set-TextValue 49 4 '0.4064'
Set-TextValue 49 5 '  +1.22%  '
Set-TextValue 50 4 '9.185'
Set-TextValue 50 5 '  +0.22%  '
Set-TextValue 51 4 '0.05797'
Set-TextValue 51 5 '  +0.81%  '
